# add Save column in s_vals sheets
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing header cell (G1) onto the new
# header cell (H1) so it keeps the same bold/border/center style, then set
# its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# New data cell for row 2: Save = 1
$ws.Range("H2").Value = 1
